# Generate Report for Handoff
#
# Refresh the "Latest Handoff Datetime" column (column D) for every file row
# whose handoff batch just completed as part of this report run (the rows
# with status "Handback transform failed" or "Ready for handoff"), on both
# the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$zhcnNewDatetime = "2016-03-08 20:29:26"
$dedeNewDatetime = "2016-03-08 20:29:35"

# Rows 7, 10-16 are the files included in this handoff batch.
$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("D" + $r).Value = $zhcnNewDatetime
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("D" + $r).Value = $dedeNewDatetime
}
